$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "precio_promocion" values in column C (rows 2-21),
# while keeping the existing cell formatting/style.
$ws.Range("C2:C21").ClearContents()

# Move the active selection to C2 (matches the saved selection state).
$ws.Range("C2").Select()
